# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# cbb64caf-... row (row 3) on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("zh-cn")
$ws1.Range("E3").Value = "2016-03-20 14:40:39"
$ws1.Range("H3").Value = "2016-03-20 14:41:00"

$ws2 = $wb.Worksheets.Item("de-de")
$ws2.Range("E3").Value = "2016-03-20 14:40:42"
$ws2.Range("H3").Value = "2016-03-20 14:41:06"
